# Update the en/ba translation sheet: re-order several rows of column A
# (the "en" source strings) to keep the list alphabetically sorted,
# add a few new strings (marked "new" in column C) and push the four
# strings that no longer fit the current ordering to the bottom of the
# sheet (marked "deleted" in column C), growing the sheet from 181 to
# 185 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 11-50: two obsolete rows removed from the top of this block
# (they reappear at the bottom, see rows 182-183) and two new rows
# inserted ("Contains names of organisms..." and "Download Lab Log...").
$ws.Cells.Item(11, 1).Value = 'ACORN Participating Countries'

$ws.Cells.Item(12, 1).Value = 'All ''orgname'' are provided.'

$ws.Cells.Item(13, 1).Value = 'All ''patid'' are provided.'

$ws.Cells.Item(14, 1).Value = 'All ''specdate'' are provided.'

$ws.Cells.Item(15, 1).Value = 'All ''specdate'' are today or before today.'

$ws.Cells.Item(16, 1).Value = 'All ''specgroup'' are provided.'

$ws.Cells.Item(17, 1).Value = 'All ''specid'' are provided.'

$ws.Cells.Item(18, 1).Value = 'All dates of enrolment for HAI patients have a matching date in the HAI survey dataset'

$ws.Cells.Item(19, 1).Value = 'All Other Organisms'

$ws.Cells.Item(20, 1).Value = 'All valid records have an ACORN ID.'

$ws.Cells.Item(21, 1).Value = 'AMR'

$ws.Cells.Item(22, 1).Value = 'and generate enrolment log.'

$ws.Cells.Item(23, 1).Value = 'Attempting to connect.'

$ws.Cells.Item(24, 1).Value = 'Blood culture collected within 24 hours of admission (CAI) / symptom onset (HAI)'

$ws.Cells.Item(25, 1).Value = 'Blood Culture Contaminants'

$ws.Cells.Item(26, 1).Value = 'Bloodstream Infection (BSI)'

$ws.Cells.Item(27, 1).Value = 'Calculated age is consistent with ''Age Category'''

$ws.Cells.Item(28, 1).Value = 'Calculated age isn''t always consistent with ''Age Category'''

$ws.Cells.Item(29, 1).Value = 'Cancel'

$ws.Cells.Item(30, 1).Value = 'Care should be taken when interpreting rates and AMR profiles where there are small numbers of cases or bacterial isolates: point estimates may be unreliable.'

$ws.Cells.Item(31, 1).Value = 'Clinical and day-28 outcomes are consistent.'

$ws.Cells.Item(32, 1).Value = 'Clinical and day-28 outcomes aren''t consistent for some dead patients.'

$ws.Cells.Item(33, 1).Value = 'Clinical Outcome'

$ws.Cells.Item(34, 1).Value = 'Clinical Outcome Status:'

$ws.Cells.Item(35, 1).Value = 'Co-resistances'

$ws.Cells.Item(36, 1).Value = 'Combine Susceptible + Intermediate'

$ws.Cells.Item(37, 1).Value = 'Consider saving .acorn file on the cloud for additional security.'

$ws.Cells.Item(38, 1).Value = 'Contains names of organisms before and after mapping.'
$ws.Cells.Item(38, 3).Value = 'new'

$ws.Cells.Item(39, 1).Value = 'Couldn''t connect to server. Please check internet access.'

$ws.Cells.Item(40, 1).Value = 'Critical errors with clinical data.'

$ws.Cells.Item(41, 1).Value = 'Culture results per specimen type'

$ws.Cells.Item(42, 1).Value = 'Data Management'

$ws.Cells.Item(43, 1).Value = 'Date of Enrolment'

$ws.Cells.Item(44, 1).Value = 'Day 28'

$ws.Cells.Item(45, 1).Value = 'Day 28 Status:'

$ws.Cells.Item(46, 1).Value = 'Diagnosis at Enrolment'

$ws.Cells.Item(47, 1).Value = 'Dismiss'

$ws.Cells.Item(48, 1).Value = 'Distribution of Enrolments'

$ws.Cells.Item(49, 1).Value = 'Download Enrolment Log (.xlsx)'

$ws.Cells.Item(50, 1).Value = 'Download Lab Log (.xlsx)'
$ws.Cells.Item(50, 3).Value = 'new'

# --- Row 70: text replaced in place, old text moved to row 184.
$ws.Cells.Item(70, 1).Value = 'HAI point prevalence by '
$ws.Cells.Item(70, 3).Value = 'new'

# --- Rows 110-124: one new row inserted at the top of the block, the
# last row of the block ("Select lab data format:") falls off the end
# and reappears at row 185.
$ws.Cells.Item(110, 1).Value = 'Remove ''Not Cultured'' specimens'
$ws.Cells.Item(110, 3).Value = 'new'

$ws.Cells.Item(111, 1).Value = 'Remove blood culture contaminants from the following visualizations'

$ws.Cells.Item(112, 1).Value = 'Reset Enrolments Filters'

$ws.Cells.Item(113, 1).Value = 'Resistance to 3rd gen. Cephalosporins Over Time'

$ws.Cells.Item(114, 1).Value = 'Resistance to Carbapenems Over Time'

$ws.Cells.Item(115, 1).Value = 'Resistance to Fluoroquinolones Over Time'

$ws.Cells.Item(116, 1).Value = 'Resistance to Oxacillin Over Time'

$ws.Cells.Item(117, 1).Value = 'Resistance to Penicillin G - meningitis Over Time'

$ws.Cells.Item(118, 1).Value = 'Resistance to Penicillin G Over Time'

$ws.Cells.Item(119, 1).Value = 'Retriving data from REDCap server.'

$ws.Cells.Item(120, 1).Value = 'Save .acorn file'

$ws.Cells.Item(121, 1).Value = 'Save acorn data'

$ws.Cells.Item(122, 1).Value = 'Save on Server'

$ws.Cells.Item(123, 1).Value = 'See Breakdown by Ward'

$ws.Cells.Item(124, 1).Value = 'See by Week'

# --- New rows 182-185: the four strings displaced from their old
# positions above, appended at the end of the sheet and flagged as
# "deleted" in column C.
$ws.Cells.Item(182, 1).Value = 'ACORN data is not of the right format. Only data generated with v2.1 (or later versions) is compatible.'
$ws.Cells.Item(182, 2).Value = 'TBT'
$ws.Cells.Item(182, 3).Value = 'deleted'

$ws.Cells.Item(183, 1).Value = 'ACORN data is not of the right format. Only data generated with v2.1 is compatible.'
$ws.Cells.Item(183, 2).Value = 'TBT'
$ws.Cells.Item(183, 3).Value = 'deleted'

$ws.Cells.Item(184, 1).Value = 'HAI point prevalence by type of ward'
$ws.Cells.Item(184, 2).Value = 'TBT'
$ws.Cells.Item(184, 3).Value = 'deleted'

$ws.Cells.Item(185, 1).Value = 'Select lab data format:'
$ws.Cells.Item(185, 2).Value = 'TBT'
$ws.Cells.Item(185, 3).Value = 'deleted'
